$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 15224
$ws1.Range("F7").Value = 411
$ws1.Range("F10").Value = 15302
$ws1.Range("F12").Value = 8841
$ws1.Range("F13").Value = 356
$ws1.Range("F15").Value = 74
$ws1.Range("F16").Value = 187
$ws1.Range("F18").Value = 187
$ws1.Range("F25").Value = 1095
$ws1.Range("F27").Value = 17
$ws1.Range("F28").Value = 60
$ws1.Range("F30").Value = 33
$ws1.Range("F31").Value = 414
$ws1.Range("F35").Value = 286
$ws1.Range("F36").Value = 437
$ws1.Range("F38").Value = 5423

# Sheet "演出" (sheet2) - F column "想去人数" update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 64

# Sheet "全部类型" (sheet4) - F column "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 15224
$ws4.Range("F7").Value = 411
$ws4.Range("F10").Value = 15302
$ws4.Range("F12").Value = 8841
$ws4.Range("F13").Value = 356
$ws4.Range("F16").Value = 74
$ws4.Range("F17").Value = 187
$ws4.Range("F19").Value = 187
$ws4.Range("F26").Value = 1095
$ws4.Range("F28").Value = 17
$ws4.Range("F29").Value = 60
$ws4.Range("F31").Value = 33
$ws4.Range("F32").Value = 64
$ws4.Range("F34").Value = 414
$ws4.Range("F38").Value = 286
$ws4.Range("F39").Value = 437
$ws4.Range("F41").Value = 5423
